$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.382.63'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '1.882.63'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''0.7139'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = '''243.01'
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '''0.08045'
$ws.Range("E8").Value = '  +4.02%  '
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("E10").Value = '  +1.49%  '
$ws.Range("D11").Value = '''0.08348'
$ws.Range("E11").Value = '  -1.98%  '
$ws.Range("D12").Value = '1.886.59'
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").Value = '''5.258'
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("D14").Value = '''0.7196'
$ws.Range("E14").Value = '  +1.40%  '
$ws.Range("E15").Value = '  +3.19%  '
$ws.Range("D16").Value = '''6.329'
$ws.Range("D17").Value = '''0.000008548'
$ws.Range("E17").Value = '  +4.27%  '
$ws.Range("D18").Value = '29.389.77'
$ws.Range("E18").Value = '  +0.26%  '
$ws.Range("D19").Value = '''242.28'
$ws.Range("E19").Value = '  +0.27%  '
$ws.Range("D20").Value = '2.138.72'
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").Value = '''7.872'
$ws.Range("E23").Value = '  +0.95%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = '''0.1589'
$ws.Range("E25").Value = '  -0.89%  '
$ws.Range("D26").Value = '''163.43'
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").Value = '''9.088'
$ws.Range("E27").Value = '  +0.55%  '
$ws.Range("D28").Value = '''18.65'
$ws.Range("E28").Value = '  +0.94%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = '''4.424'
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("D31").Value = '''4.334'
$ws.Range("E31").Value = '  +0.40%  '
$ws.Range("D32").Value = '''1.200'
$ws.Range("E32").Value = '  -6.67%  '
$ws.Range("E33").Value = '  +2.47%  '
$ws.Range("E34").Value = '  +1.06%  '
$ws.Range("D35").Value = '''1.184'
$ws.Range("E35").Value = '  +0.77%  '
$ws.Range("D36").Value = '''0.7504'
$ws.Range("E36").Value = '  +0.87%  '
$ws.Range("D37").Value = '''2.698'
$ws.Range("E37").Value = '  +0.44%  '
$ws.Range("E38").Value = '  +1.36%  '
$ws.Range("D39").Value = '1.286.55'
$ws.Range("E39").Value = '  +8.72%  '
$ws.Range("D40").Value = '''2.745'
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("D41").Value = '''6.595'
$ws.Range("E41").Value = '  +3.35%  '
$ws.Range("D42").Value = '''0.9145'
$ws.Range("E42").Value = '  +2.93%  '
$ws.Range("D43").Value = '''74.63'
$ws.Range("D44").Value = '''111.87'
$ws.Range("E44").Value = '  +5.28%  '
$ws.Range("D45").Value = '''1.000'
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = '''0.00000000130'
$ws.Range("E46").Value = '  +6.96%  '
$ws.Range("D47").Value = '2.039.11'
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("D48").Value = '''1.812'
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("D49").Value = '''0.5219'
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("D50").Value = '''9.538'
$ws.Range("E50").Value = '  +1.82%  '
$ws.Range("D51").Value = '''0.4396'
$ws.Range("E51").Value = '  +1.94%  '
